# "Ready for handoff" — the b.md row has a new handoff generated (zh-cn + de-de),
# so the Overview + per-locale sheets get refreshed status/timestamp/filename/error info.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84b3607f06a861ebc61cb87bddacc620cf9f3c0b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c5a5e65af44c3c827ce1686b3dea0711e4a03f2/e2e/b.md."

# --- Overview sheet: row for b.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-30 02:39:45"

# --- zh-cn sheet: row for b.md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"

# F3 holds the literal text "False" (was "True") - use copy/paste-values so Excel
# keeps storing it as text instead of auto-coercing the look-alike word to a Boolean.
$wsZhCn.Range("F2").Copy() | Out-Null
$wsZhCn.Range("F3").PasteSpecial(-4163) | Out-Null

$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-30 02:39:41"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664

# --- de-de sheet: row for b.md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"

$wsDeDe.Range("F2").Copy() | Out-Null
$wsDeDe.Range("F3").PasteSpecial(-4163) | Out-Null

$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-30 02:39:45"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1").EntireColumn.ColumnWidth = 39.166666666666664
